$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HEATING")
$ws.Activate()

# Give the brand-new row 7 the same per-cell formatting as row 6 (its
# cells don't exist yet, so they'd otherwise get no style at all).
$ws.Range("A6:I6").Copy()
$ws.Range("A7:I7").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Shift existing rows 4-6 down to rows 5-7 (bottom-up, to avoid clobbering
# source data before it is copied). Value2 (unlike Value in this engine)
# correctly round-trips both read and write, and direct cell writes keep
# each cell's existing style index.
$ws.Range("A7").Value2 = $ws.Range("A6").Value2
$ws.Range("B7").Value2 = $ws.Range("B6").Value2
$ws.Range("C7").Value2 = $ws.Range("C6").Value2
$ws.Range("D7").Value2 = $ws.Range("D6").Value2
$ws.Range("E7").Value2 = $ws.Range("E6").Value2
$ws.Range("F7").Value2 = $ws.Range("F6").Value2
$ws.Range("G7").Value2 = $ws.Range("G6").Value2
$ws.Range("H7").Value2 = $ws.Range("H6").Value2

$ws.Range("A6").Value2 = $ws.Range("A5").Value2
$ws.Range("B6").Value2 = $ws.Range("B5").Value2
$ws.Range("C6").Value2 = $ws.Range("C5").Value2
$ws.Range("D6").Value2 = $ws.Range("D5").Value2
$ws.Range("E6").Value2 = $ws.Range("E5").Value2
$ws.Range("F6").Value2 = $ws.Range("F5").Value2
$ws.Range("G6").Value2 = $ws.Range("G5").Value2
$ws.Range("H6").Value2 = $ws.Range("H5").Value2

$ws.Range("A5").Value2 = $ws.Range("A4").Value2
$ws.Range("B5").Value2 = $ws.Range("B4").Value2
$ws.Range("C5").Value2 = $ws.Range("C4").Value2
$ws.Range("D5").Value2 = $ws.Range("D4").Value2
$ws.Range("E5").Value2 = $ws.Range("E4").Value2
$ws.Range("F5").Value2 = $ws.Range("F4").Value2
$ws.Range("G5").Value2 = $ws.Range("G4").Value2
$ws.Range("H5").Value2 = $ws.Range("H4").Value2

# Populate row 4 with the new "natural gas-fired boiler" entry.
$ws.Range("A4").Value2 = "natural gas-fired boiler"
$ws.Range("B4").Value2 = "T3"
$ws.Range("C4").Value2 = "NG"
$ws.Range("D4").Value2 = 0.8
$ws.Range("E4").Value2 = 1.403
$ws.Range("F4").Value2 = 0.1
$ws.Range("G4").Value2 = 0.22
$ws.Range("H4").Value2 = "KBOB 2019, costs in USD-2015"

# Update the active selection to match the committed state.
$ws.Range("A11").Select()
